{"js": "// Circle Language Spec Plan: More uniform mentioning of months in the\n// sub-project names.\n//\n// The document title paragraph (Heading2) reads:\n//   \"Circle Language Spec Plan,\" <break> \"Classes & Relations Specs,\" <break> \"Project Summary\"\n// The sub-project name gets its month prefixed, becoming:\n//   \"Circle Language Spec Plan,\" <break> \"2008-05 Classes & Relations Specs,\" <break> \"Project Summary\"\n//\n// (The document also carries a hidden \"_GoBack\" bookmark - Word's marker\n// for \"the last place I typed\" - which ends up sitting right where the new\n// text was typed, i.e. immediately before \"Classes\".)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The title is the first paragraph in the document.\nconst titleParagraph = paragraphs.items[0];\n\n// Find \"Classes\" inside that paragraph only, so we don't touch the other\n// seven occurrences of the word elsewhere in the document.\nconst matches = titleParagraph.search(\"Classes\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nconst classesRange = matches.items[0];\n\n// Type the new sub-project month prefix right before \"Classes\".\nclassesRange.insertText(\"2008-05 \", Word.InsertLocation.before);\nawait context.sync();\n\n// Re-anchor Word's \"_GoBack\" last-edit bookmark to the freshly typed spot,\n// immediately before \"Classes\" (mirroring Word's own behaviour of moving\n// _GoBack to track the most recent edit location).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst refreshed = titleParagraph.search(\"Classes\", { matchCase: true, matchWholeWord: false });\nrefreshed.load(\"items\");\nawait context.sync();\n\nconst newEditSpot = refreshed.items[0].getRange(Word.RangeLocation.start);\nnewEditSpot.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Circle Language Spec Plan: More uniform mentioning of months in the\n# sub-project names.\n#\n# The document title paragraph (Heading2) reads:\n#   \"Circle Language Spec Plan,\" <break> \"Classes & Relations Specs,\" <break> \"Project Summary\"\n# The sub-project name gets its month prefixed, becoming:\n#   \"Circle Language Spec Plan,\" <break> \"2008-05 Classes & Relations Specs,\" <break> \"Project Summary\"\n#\n# (The document also carries a hidden \"_GoBack\" bookmark - Word's marker\n# for \"the last place I typed\" - which ends up sitting right where the new\n# text was typed, i.e. immediately before \"Classes\".)\n\n$d = $word.ActiveDocument\n\n# Find \"Classes\" - it is the start of the sub-project name in the title\n# (first paragraph / Heading2). Using Find on the title paragraph's range\n# keeps this targeted at that one occurrence.\n$titleRange = $d.Paragraphs(1).Range\n$titleRange.Find.MatchCase = $true\n$titleRange.Find.Execute(\"Classes\") | Out-Null\n\n# Type the new sub-project month prefix right before \"Classes\".\n$insertionPoint = $d.Range($titleRange.Start, $titleRange.Start)\n$insertionPoint.InsertBefore(\"2008-05 \")\n\n# Re-anchor Word's \"_GoBack\" last-edit bookmark to the freshly typed spot,\n# immediately before \"Classes\" (mirroring Word's own behaviour of moving\n# _GoBack to track the most recent edit location).\n$d.Bookmarks(\"_GoBack\").Delete()\n\n$retitleRange = $d.Paragraphs(1).Range\n$retitleRange.Find.MatchCase = $true\n$retitleRange.Find.Execute(\"Classes\") | Out-Null\n$newEditSpot = $d.Range($retitleRange.Start, $retitleRange.Start)\n$d.Bookmarks.Add(\"_GoBack\", $newEditSpot) | Out-Null\n"}
